# Update column F (dSF) values for specific rows to repulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -7
    "F4"  = 3
    "F7"  = -4
    "F19" = -6
    "F20" = -4
    "F21" = 2
    "F24" = -1
    "F27" = 2
    "F34" = 7
    "F38" = -1
    "F41" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
